$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values remain plain text (matching original inline-string formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values
$ws.Range('D2').Value = '68.398.47'
$ws.Range('E2').Value = '  -4.43%  '
$ws.Range('D3').Value = '3.712.30'
$ws.Range('E3').Value = '  -4.39%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '584.78'
$ws.Range('E5').Value = '  -2.34%  '
$ws.Range('D6').Value = '180.65'
$ws.Range('E6').Value = '  +7.36%  '
$ws.Range('D7').Value = '3.705.74'
$ws.Range('E7').Value = '  -4.48%  '
$ws.Range('D8').Value = '0.629'
$ws.Range('E8').Value = '  -6.53%  '
$ws.Range('D9').Value = '0.997'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D10').Value = '0.715'
$ws.Range('E10').Value = '  -5.86%  '
$ws.Range('E11').Value = '  -7.46%  '
$ws.Range('D12').Value = '54.07'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').Value = '0.0000292'
$ws.Range('E13').Value = '  -9.21%  '
$ws.Range('D14').Value = '10.47'
$ws.Range('E14').Value = '  -8.21%  '
$ws.Range('D15').Value = '4.204.05'
$ws.Range('E15').Value = '  -6.69%  '
$ws.Range('D16').Value = '3.705.20'
$ws.Range('E16').Value = '  -4.94%  '
$ws.Range('D17').Value = '19.55'
$ws.Range('E17').Value = '  -6.61%  '
$ws.Range('E18').Value = '  -2.68%  '
$ws.Range('D19').Value = '12.87'
$ws.Range('E19').Value = '  -7.54%  '
$ws.Range('D20').Value = '1.13'
$ws.Range('E20').Value = '  -7.64%  '
$ws.Range('D21').Value = '67.999.89'
$ws.Range('E21').Value = '  -4.74%  '
$ws.Range('D22').Value = '409.62'
$ws.Range('E22').Value = '  -6.02%  '
$ws.Range('E23').Value = '  -5.40%  '
$ws.Range('D24').Value = '88.68'
$ws.Range('E24').Value = '  -6.02%  '
$ws.Range('E25').Value = '  -8.62%  '
$ws.Range('D26').Value = '12.85'
$ws.Range('E26').Value = '  -7.37%  '
$ws.Range('D27').Value = '11.09'
$ws.Range('E27').Value = '  +0.64%  '
$ws.Range('D28').Value = '3.88'
$ws.Range('E28').Value = '  -6.77%  '
$ws.Range('D29').Value = '6.06'
$ws.Range('E29').Value = '  +1.96%  '
$ws.Range('D30').Value = '9.52'
$ws.Range('E30').Value = '  -6.93%  '
$ws.Range('D31').Value = '32.66'
$ws.Range('E31').Value = '  -6.94%  '
$ws.Range('D32').Value = '7.51'
$ws.Range('E32').Value = '  -6.22%  '
$ws.Range('D33').Value = '12.57'
$ws.Range('E33').Value = '  -8.40%  '
$ws.Range('E34').Value = '  -6.99%  '
$ws.Range('D35').Value = '65.31'
$ws.Range('E35').Value = '  -4.58%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = '43.49'
$ws.Range('E36').Value = '  -16.65%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').Value = '601.79'
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('D38').Value = '0.0₃0901'
$ws.Range('E38').Value = '  -9.99%  '
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').Value = '0.400'
$ws.Range('E40').Value = '  -5.55%  '
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  -4.12%  '
$ws.Range('D43').Value = '2.79'
$ws.Range('E43').Value = '  +4.17%  '
$ws.Range('D44').Value = '3.00'
$ws.Range('E44').Value = '  -9.33%  '
$ws.Range('D45').Value = '2.96'
$ws.Range('E45').Value = '  -9.88%  '
$ws.Range('D46').Value = '0.0436'
$ws.Range('E46').Value = '  -7.23%  '
$ws.Range('D47').Value = '9.28'
$ws.Range('E47').Value = '  -9.10%  '
$ws.Range('D48').Value = '2.802.03'
$ws.Range('E48').Value = '  -1.88%  '
$ws.Range('D49').Value = '0.134'
$ws.Range('E49').Value = '  -7.44%  '
$ws.Range('E50').Value = '  -3.95%  '
$ws.Range('D51').Value = '3.14'
$ws.Range('E51').Value = '  -5.79%  '
